# Auto-generated edit script applying the crypto price/volume refresh diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.136.96'
$ws.Range("E2").Value = '  -0.69%  '
$ws.Range("D3").Value = '3.323.97'
$ws.Range("E3").Value = '  -1.35%  '
$ws.Range("E4").Value = '  -0.11%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '551.49'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.92%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '172.87'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.93%  '
$ws.Range("E7").Value = '  +1.25%  '
$ws.Range("E8").Value = '  +0.10%  '
$ws.Range("D9").Value = '3.315.54'
$ws.Range("E9").Value = '  -1.41%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.171'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +5.70%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.638'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.58%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '53.25'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.90%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000277'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.80%  '
$ws.Range("E14").Value = '  -0.21%  '
$ws.Range("D15").Value = '3.855.95'
$ws.Range("E15").Value = '  -1.49%  '
$ws.Range("E16").Value = '  +2.14%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '18.07'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.81%  '
$ws.Range("D18").Value = '3.322.95'
$ws.Range("E18").Value = '  -1.74%  '
$ws.Range("D19").Value = '64.299.47'
$ws.Range("E19").Value = '  -0.40%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.69'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.34%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.981'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.07%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '447.38'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +3.48%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.96'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.51%  '
$ws.Range("E24").Value = '  -2.10%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '86.65'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.85%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '13.72'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +3.44%  '
$ws.Range("E27").Value = '  +0.92%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.60'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.85%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.56'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -2.35%  '
$ws.Range("E30").Value = '  +3.37%  '
$ws.Range("E31").Value = '  -2.80%  '
$ws.Range("B32").Value = 'Cosmos'
$ws.Range("C32").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '11.35'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.25%  '
$ws.Range("B33").Value = 'OKB'
$ws.Range("C33").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '62.19'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +6.48%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '569.33'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.43%  '
$ws.Range("E35").Value = '  -1.84%  '
$ws.Range("E36").Value = '  +0.04%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.55'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.48%  '
$ws.Range("E38").Value = '  -1.30%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '35.14'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.77%  '
$ws.Range("E40").Value = '  -1.09%  '
$ws.Range("D41").Value = '0.0₃0728'
$ws.Range("E41").Value = '  -4.08%  '
$ws.Range("D42").Value = '3.056.46'
$ws.Range("E42").Value = '  -1.89%  '
$ws.Range("E43").Value = '  +0.45%  '
$ws.Range("E44").Value = '  -4.25%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.16'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -3.23%  '
$ws.Range("E46").Value = '  +2.42%  '
$ws.Range("E47").Value = '  -1.40%  '
$ws.Range("B48").Value = 'Monero'
$ws.Range("C48").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '142.43'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +5.72%  '
$ws.Range("B49").Value = 'FirstDigitalUSD'
$ws.Range("C49").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.00'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.09%  '
$ws.Range("E50").Value = '  -3.18%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '8.16'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.24%  '

Write-Host "Applied 90 cell updates"